$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "model_4_5_23"
$ws.Range("B2").Value = 0.3622594757142954
$ws.Range("C2").Value = -0.5125150414117985
$ws.Range("D2").Value = -1.373887891416268
$ws.Range("E2").Value = -0.6145112206551215
$ws.Range("F2").Value = 0.7057908773422241
$ws.Range("G2").Value = 3.211472988128662
$ws.Range("H2").Value = 0.7844958305358887
$ws.Range("I2").Value = 2.069365501403809

$ws.Range("A3").Value = "model_4_5_22"
$ws.Range("B3").Value = 0.3625594857168539
$ws.Range("C3").Value = -0.5090652396989142
$ws.Range("D3").Value = -1.385674718229159
$ws.Range("E3").Value = -0.6129216683398833
$ws.Range("F3").Value = 0.7054587602615356
$ws.Range("G3").Value = 3.204148054122925
$ws.Range("H3").Value = 0.7883910536766052
$ws.Range("I3").Value = 2.067327976226807

$ws.Range("A4").Value = "model_4_5_24"
$ws.Range("B4").Value = 0.3644254426447037
$ws.Range("C4").Value = -0.5062035510141583
$ws.Range("D4").Value = -1.373272822653448
$ws.Range("E4").Value = -0.6089017910039591
$ws.Range("F4").Value = 0.7033937573432922
$ws.Range("G4").Value = 3.198071956634521
$ws.Range("H4").Value = 0.7842925786972046
$ws.Range("I4").Value = 2.062175750732422

$ws.Range("A5").Value = "model_4_5_21"
$ws.Range("B5").Value = 0.3836817245586538
$ws.Range("C5").Value = -0.4417067825090832
$ws.Range("D5").Value = -1.368373021065424
$ws.Range("E5").Value = -0.551743305084841
$ws.Range("F5").Value = 0.6820827126502991
$ws.Range("G5").Value = 3.06112813949585
$ws.Range("H5").Value = 0.7826733589172363
$ws.Range("I5").Value = 1.988914012908936

$ws.Range("A6").Value = "model_4_5_20"
$ws.Range("B6").Value = 0.3930389955100585
$ws.Range("C6").Value = -0.4156570195912068
$ws.Range("D6").Value = -1.333441921891596
$ws.Range("E6").Value = -0.524658778886735
$ws.Range("F6").Value = 0.6717268824577332
$ws.Range("G6").Value = 3.005817413330078
$ws.Range("H6").Value = 0.7711296677589417
$ws.Range("I6").Value = 1.954198956489563

$ws.Range("A7").Value = "model_4_5_18"
$ws.Range("B7").Value = 0.3960253337480781
$ws.Range("C7").Value = -0.4117614025005183
$ws.Range("D7").Value = -1.273118650116455
$ws.Range("E7").Value = -0.5139175776088214
$ws.Range("F7").Value = 0.6684218645095825
$ws.Range("G7").Value = 2.997546195983887
$ws.Range("H7").Value = 0.7511947154998779
$ws.Range("I7").Value = 1.940431714057922

$ws.Range("A8").Value = "model_4_5_19"
$ws.Range("B8").Value = 0.3969398033090307
$ws.Range("C8").Value = -0.4110427281712548
$ws.Range("D8").Value = -1.261707348131676
$ws.Range("E8").Value = -0.5119088092061324
$ws.Range("F8").Value = 0.6674098968505859
$ws.Range("G8").Value = 2.996020317077637
$ws.Range("H8").Value = 0.7474236488342285
$ws.Range("I8").Value = 1.937856793403625

$ws.Range("A9").Value = "model_4_5_17"
$ws.Range("B9").Value = 0.3971325689249711
$ws.Range("C9").Value = -0.4090234179243795
$ws.Range("D9").Value = -1.264029291709566
$ws.Range("E9").Value = -0.5104196220843409
$ws.Range("F9").Value = 0.6671966314315796
$ws.Range("G9").Value = 2.991732597351074
$ws.Range("H9").Value = 0.7481909394264221
$ws.Range("I9").Value = 1.935948133468628

$ws.Range("A10").Value = "model_4_5_16"
$ws.Range("B10").Value = 0.4093717715288218
$ws.Range("C10").Value = -0.3731425435499078
$ws.Range("D10").Value = -1.234456459356049
$ws.Range("E10").Value = -0.4753637661118426
$ws.Range("F10").Value = 0.6536513566970825
$ws.Range("G10").Value = 2.915547847747803
$ws.Range("H10").Value = 0.7384180426597595
$ws.Range("I10").Value = 1.891016364097595

$ws.Range("A11").Value = "model_4_5_15"
$ws.Range("B11").Value = 0.4118840472385233
$ws.Range("C11").Value = -0.367791198423268
$ws.Range("D11").Value = -1.20850885561527
$ws.Range("E11").Value = -0.4675223896584073
$ws.Range("F11").Value = 0.6508709788322449
$ws.Range("G11").Value = 2.90418553352356
$ws.Range("H11").Value = 0.7298431992530823
$ws.Range("I11").Value = 1.880965709686279

$ws.Range("A12").Value = "model_4_5_14"
$ws.Range("B12").Value = 0.4251328212203263
$ws.Range("C12").Value = -0.3314779280929561
$ws.Range("D12").Value = -1.154618868641784
$ws.Range("E12").Value = -0.4291368301695995
$ws.Range("F12").Value = 0.6362085342407227
$ws.Range("G12").Value = 2.827082633972168
$ws.Range("H12").Value = 0.712034285068512
$ws.Range("I12").Value = 1.83176577091217

$ws.Range("A13").Value = "model_4_5_13"
$ws.Range("B13").Value = 0.4267130217933884
$ws.Range("C13").Value = -0.3276738551446523
$ws.Range("D13").Value = -1.141011950262224
$ws.Range("E13").Value = -0.4241433730646018
$ws.Range("F13").Value = 0.6344597935676575
$ws.Range("G13").Value = 2.819005489349365
$ws.Range("H13").Value = 0.7075375914573669
$ws.Range("I13").Value = 1.825365424156189

$ws.Range("A14").Value = "model_4_5_12"
$ws.Range("B14").Value = 0.4347922939649401
$ws.Range("C14").Value = -0.3053120743525555
$ws.Range("D14").Value = -1.111420172188845
$ws.Range("E14").Value = -0.4009479417836765
$ws.Range("F14").Value = 0.6255183815956116
$ws.Range("G14").Value = 2.771525621414185
$ws.Range("H14").Value = 0.6977584362030029
$ws.Range("I14").Value = 1.795635342597961

$ws.Range("A15").Value = "model_4_5_10"
$ws.Range("B15").Value = 0.4360263234811244
$ws.Range("C15").Value = -0.3071710021624132
$ws.Range("D15").Value = -1.063889938869722
$ws.Range("E15").Value = -0.3968110522609007
$ws.Range("F15").Value = 0.624152660369873
$ws.Range("G15").Value = 2.77547287940979
$ws.Range("H15").Value = 0.6820511221885681
$ws.Range("I15").Value = 1.790332794189453

$ws.Range("A16").Value = "model_4_5_11"
$ws.Range("B16").Value = 0.4370633115721536
$ws.Range("C16").Value = -0.2988399107545177
$ws.Range("D16").Value = -1.102274210818058
$ws.Range("E16").Value = -0.3941620985143819
$ws.Range("F16").Value = 0.6230049729347229
$ws.Range("G16").Value = 2.75778341293335
$ws.Range("H16").Value = 0.6947360038757324
$ws.Range("I16").Value = 1.786937832832336

$ws.Range("A17").Value = "model_4_5_9"
$ws.Range("B17").Value = 0.43926929678847
$ws.Range("C17").Value = -0.3049412458345044
$ws.Range("D17").Value = -0.9688307266766614
$ws.Range("E17").Value = -0.3833158014807394
$ws.Range("F17").Value = 0.6205636262893677
$ws.Range("G17").Value = 2.770738363265991
$ws.Range("H17").Value = 0.6506370902061462
$ws.Range("I17").Value = 1.773035645484924

$ws.Range("A18").Value = "model_4_5_8"
$ws.Range("B18").Value = 0.4677507136954818
$ws.Range("C18").Value = -0.2412433008828401
$ws.Range("D18").Value = -0.7392018138249805
$ws.Range("E18").Value = -0.2995972388187813
$ws.Range("F18").Value = 0.589043140411377
$ws.Range("G18").Value = 2.635490655899048
$ws.Range("H18").Value = 0.5747519135475159
$ws.Range("I18").Value = 1.665731191635132

$ws.Range("A19").Value = "model_4_5_7"
$ws.Range("B19").Value = 0.4821529743140259
$ws.Range("C19").Value = -0.2077199448238625
$ws.Range("D19").Value = -0.597169897333103
$ws.Range("E19").Value = -0.2529643080273858
$ws.Range("F19").Value = 0.5731040239334106
$ws.Range("G19").Value = 2.564311504364014
$ws.Range("H19").Value = 0.5278148055076599
$ws.Range("I19").Value = 1.605960369110107

$ws.Range("A20").Value = "model_4_5_6"
$ws.Range("B20").Value = 0.5241549309862923
$ws.Range("C20").Value = -0.08090229106970859
$ws.Range("D20").Value = -0.519787050731785
$ws.Range("E20").Value = -0.1323484786743252
$ws.Range("F20").Value = 0.5266203284263611
$ws.Range("G20").Value = 2.2950439453125
$ws.Range("H20").Value = 0.5022421479225159
$ws.Range("I20").Value = 1.451363563537598

$ws.Range("A21").Value = "model_4_5_5"
$ws.Range("B21").Value = 0.5494752084743664
$ws.Range("C21").Value = -0.007037461340150264
$ws.Range("D21").Value = -0.4048556152670135
$ws.Range("E21").Value = -0.05363801226765519
$ws.Range("F21").Value = 0.4985981583595276
$ws.Range("G21").Value = 2.138209104537964
$ws.Range("H21").Value = 0.4642608761787415
$ws.Range("I21").Value = 1.350478053092957

$ws.Range("A22").Value = "model_4_5_4"
$ws.Range("B22").Value = 0.6270015803126774
$ws.Range("C22").Value = 0.2440366913989859
$ws.Range("D22").Value = -0.245966230622173
$ws.Range("E22").Value = 0.1858401027742906
$ws.Range("F22").Value = 0.4127993583679199
$ws.Range("G22").Value = 1.605111718177795
$ws.Range("H22").Value = 0.4117529392242432
$ws.Range("I22").Value = 1.043532133102417

$ws.Range("A23").Value = "model_4_5_0"
$ws.Range("B23").Value = 0.7273286034814863
$ws.Range("C23").Value = 0.9149962119211804
$ws.Range("D23").Value = 0.6041102709755928
$ws.Range("E23").Value = 0.8774170518561198
$ws.Range("F23").Value = 0.3017669022083282
$ws.Range("G23").Value = 0.180485725402832
$ws.Range("H23").Value = 0.1308292001485825
$ws.Range("I23").Value = 0.1571180820465088

$ws.Range("A24").Value = "model_4_5_3"
$ws.Range("B24").Value = 0.7396013200687955
$ws.Range("C24").Value = 0.661403856422475
$ws.Range("D24").Value = 0.006353665361202077
$ws.Range("E24").Value = 0.5824884391728402
$ws.Range("F24").Value = 0.2881845831871033
$ws.Range("G24").Value = 0.7189299464225769
$ws.Range("H24").Value = 0.3283690810203552
$ws.Range("I24").Value = 0.5351365804672241

$ws.Range("A25").Value = "model_4_5_2"
$ws.Range("B25").Value = 0.7502231566978884
$ws.Range("C25").Value = 0.6660524697908434
$ws.Range("D25").Value = 0.2248658690512504
$ws.Range("E25").Value = 0.6130780084102597
$ws.Range("F25").Value = 0.2764293551445007
$ws.Range("G25").Value = 0.7090596556663513
$ws.Range("H25").Value = 0.2561576068401337
$ws.Range("I25").Value = 0.495929092168808

$ws.Range("A26").Value = "model_4_5_1"
$ws.Range("B26").Value = 0.7579143690161017
$ws.Range("C26").Value = 0.6968620497532007
$ws.Range("D26").Value = 0.3324950351417244
$ws.Range("E26").Value = 0.6531566534980504
$ws.Range("F26").Value = 0.2679174542427063
$ws.Range("G26").Value = 0.6436427235603333
$ws.Range("H26").Value = 0.2205895632505417
$ws.Range("I26").Value = 0.4445590972900391
